# Update the cached "today" date shown by the datetimeFigureOut field
# placeholders on the slide master and every slide layout (4/21/2021 -> 6/24/2025).
$p = $ppt.ActivePresentation
$newDate = "6/24/2025"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# Rename the "Refinitiv Real-Time / Optimized (AWS)" box on slide 1 to
# "Real-Time / Optimized (RTO)".
$slide1 = $p.Slides.Item(1)
$rect22 = $slide1.Shapes.Item("Rectangle 22")
$tr = $rect22.TextFrame.TextRange
$tr.Paragraphs(1).Runs(1).Text = "Real-Time"
$tr.Paragraphs(2).Runs(1).Text = "Optimized (RTO)"
